# Trade #2 closed at 2026-02-16 22:56:31 - base_strategy UP +0.000%
# Append the new trade row (row 3) to both the "All Trades" and
# "base_strategy" worksheets.

$wb = $excel.ActiveWorkbook

$sheetNames = @("All Trades", "base_strategy")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)

    $ws.Cells.Item(3, 1).Value = 2                  # A3 Trade #
    $ws.Cells.Item(3, 2).Value = "'2026-02-16"      # B3 Date (force text, not a date serial)
    $ws.Cells.Item(3, 3).Value = "'22:56:31"        # C3 Time (force text, not a time serial)
    $ws.Cells.Item(3, 4).Value = "base_strategy"    # D3 Strategy
    $ws.Cells.Item(3, 5).Value = "UP"               # E3 Side
    $ws.Cells.Item(3, 6).Value = 0.5                # F3 Entry Price
    $ws.Cells.Item(3, 7).Value = "'"                # G3 Exit Price (empty text, still open)
    $ws.Cells.Item(3, 8).Value = "OPEN"             # H3 Status
    $ws.Cells.Item(3, 9).Value = 0                  # I3 P&L %
    $ws.Cells.Item(3, 10).Value = 0                 # J3 P&L $
    $ws.Cells.Item(3, 11).Value = 100               # K3 Capital After
    $ws.Cells.Item(3, 12).Value = 0                 # L3 Entry Slippage (bps)
    $ws.Cells.Item(3, 13).Value = 0                 # M3 Exit Slippage (bps)
    $ws.Cells.Item(3, 14).Value = 0.6               # N3 Confidence
    $ws.Cells.Item(3, 15).Value = "Normal spread capture: 19600 bps"  # O3 Entry Reason
    $ws.Cells.Item(3, 16).Value = "'"               # P3 Exit Reason (empty text, still open)
    $ws.Cells.Item(3, 17).Value = 0                 # Q3 Duration (min)
}
